$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B5 from text "1" to a numeric value 1
$ws.Range("B5").Value = 1

# Add new row 6
$ws.Range("A6").Value = "Ying Tang"
$ws.Range("B6").Value = "3"
$ws.Range("C6").Value = "无"
$ws.Range("D6").Value = "SUG"
$ws.Range("E6").Value = "EXP"
$ws.Range("F6").Value = "53dcf950-aee9-43ba-bb93-9e7c5cd5833d"
$ws.Range("G6").Value = "By5SY2gA-_annotated.xlsx"
$ws.Range("H6").Value = "For instance, what about averaging WordNet path-based distance metrics and distance in word embedding space (for word similarity), and other ways of applying the affect data to email tone prediction?"
